# Rename the "edge" sheet's header cells from Orig/Dest to from/to,
# and move the active selection to B2 (matching the author's re-upload).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("edge")

$ws.Range("A1").Value = "from"
$ws.Range("B1").Value = "to"

[void]$ws.Activate()
[void]$ws.Range("B2").Select()
